# Add team record (Wins / Losses / Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, right after the existing "Unnamed: 28" column (AC).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the formatting (bold, centered, bordered) used by the other header cells.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in the team's record for every data row (2-46): 86 wins, 76 losses, 0 ties.
$lastRow = 46
$ws.Range("AD2:AD" + $lastRow).Value = 86
$ws.Range("AE2:AE" + $lastRow).Value = 76
$ws.Range("AF2:AF" + $lastRow).Value = 0
